# issue359/test_01.xlsx
# refactor: update demo data and csv file for testing
#
# - Replace the Chinese column headers with English equivalents.
# - Resize columns A and B (custom widths).
# - Move the active cell selection from E9 to B14.
# - Adjust the window height.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text: Chinese -> English ------------------------------------
$ws.Range("A1").Value = "StringTitle"
$ws.Range("B1").Value = "DateTitle"
$ws.Range("C1").Value = "DoubleTitle"

# --- Column widths for columns A and B -----------------------------------
# Target stored widths (in characters) are ~14.41 and ~16.66; ColumnWidth
# is quantized internally to the default-font step, so these inputs land
# on the closest reachable stored widths.
$ws.Columns.Item(1).ColumnWidth = 13.7
$ws.Columns.Item(2).ColumnWidth = 16

# --- Window height (bookViews/workbookView@windowHeight is in twips;
#     Window.Height is in points, so 27520 twips -> 1376 pt) -------------
$wb.Windows.Item(1).Height = 1376

# --- Active selection moves from E9 to B14 --------------------------------
$ws.Range("B14").Select() | Out-Null
